# Update data and import_data for Fall 2019 dump
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Append the new Fall 2019 data row (row 34), re-using the existing
#    row's cell formats so no extra number formats/styles get minted.
$ws.Range("A33:G33").Copy()
$ws.Range("A34").PasteSpecial(-4122)

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 43747.9993055556
$ws.Range("C34").Value = 43745.9993055556
$ws.Range("D34").Value = "2019 October Lobbyist Report"
$ws.Range("E34").Value = 1
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 43592

# 2. Rename sheet to reflect the NMInDepth prefix used for the Fall 2019 dump
$ws.Name = "NMInDepth_Cam_FilingPeriodLobby"

# 3. Update the defined name to match the new sheet name/range
$wb.Names("Cam_FilingPeriodLobbyist").Delete()
$wb.Names.Add("NMInDepth_Cam_FilingPeriodLobbyist", "='NMInDepth_Cam_FilingPeriodLobby'!`$A`$1:`$H`$34")
